$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9167267680168152
$ws.Range("B1").Value = 1.790721416473389
$ws.Range("C1").Value = 4.097220897674561
$ws.Range("D1").Value = 3.523437976837158
$ws.Range("E1").Value = 1.491192221641541
